# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (F) / "最低票价" (G) figures across the three sheets
# that carry this event data: 展览, 演出, 全部类型.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 3045
$ws.Range("G6").Value = 88
$ws.Range("F10").Value = 779
$ws.Range("F11").Value = 357
$ws.Range("F12").Value = 4539
$ws.Range("F13").Value = 4539
$ws.Range("F14").Value = 104
$ws.Range("F16").Value = 140
$ws.Range("F19").Value = 89
$ws.Range("F20").Value = 7097
$ws.Range("F24").Value = 491
$ws.Range("F25").Value = 1292
$ws.Range("G25").Value = 70
$ws.Range("F26").Value = 6271
$ws.Range("F27").Value = 1664
$ws.Range("F29").Value = 1978
$ws.Range("F30").Value = 6075
$ws.Range("F34").Value = 90
$ws.Range("F36").Value = 6196
$ws.Range("F38").Value = 195
$ws.Range("F40").Value = 19
$ws.Range("F41").Value = 13
$ws.Range("F42").Value = 2432
$ws.Range("F45").Value = 1014
$ws.Range("F47").Value = 382
$ws.Range("F48").Value = 2095
$ws.Range("F49").Value = 27

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 34

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 3045
$ws.Range("G6").Value = 88
$ws.Range("F11").Value = 357
$ws.Range("F12").Value = 4539
$ws.Range("F13").Value = 4539
$ws.Range("F14").Value = 104
$ws.Range("F16").Value = 140
$ws.Range("F19").Value = 89
$ws.Range("F20").Value = 7097
$ws.Range("F23").Value = 491
$ws.Range("F24").Value = 1292
$ws.Range("G24").Value = 70
$ws.Range("F26").Value = 6271
$ws.Range("F27").Value = 1664
$ws.Range("F28").Value = 1978
$ws.Range("F31").Value = 6075
$ws.Range("F36").Value = 90
$ws.Range("F38").Value = 6196
$ws.Range("F40").Value = 195
$ws.Range("F42").Value = 13
$ws.Range("F44").Value = 2432
$ws.Range("F46").Value = 1014
$ws.Range("F48").Value = 382
$ws.Range("F49").Value = 2095
$ws.Range("F50").Value = 27
